# Infrastructures de Santé workbook - corrections/clean-up pass
# "Ajout des tableaux corrige et affichage des colonnes vides"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Centres de Santé"
$ws2 = $wb.Worksheets.Item(2)   # "Postes de Santé"

# ---------------------------------------------------------------------
# 1) Text / label corrections (order chosen to mirror the new strings
#    as they were introduced in the corrected workbook)
# ---------------------------------------------------------------------
$ws1.Range("A13").Value = "Tiris Zemmour"
$ws1.Range("A1").Value  = "Tableau : Evolution du nombre de centres de santé par Wilaya durant la période 2014-2024"
$ws2.Range("A1").Value  = "Tableau : Evolution du nombre de postes de santé par Wilaya durant la période 2014-2024"
$ws2.Range("A3").Value  = "Hodh Chargui"
$ws1.Range("A19").Value = "Source : Annuaires des Statistiques Sanitaires/MS"
$ws2.Range("A19").Value = "Source : Annuaires des Statistiques Sanitaires/MS"
$ws2.Range("A2").Value  = "Wilaya"
$ws1.Range("A4").Value  = "Hodh El Gharbi"

# ---------------------------------------------------------------------
# 2) Header row clean-up on "Postes de Santé": give the newly-labelled
#    A2 ("Wilaya") a proper header look, borrowing D2's current
#    (bold / bordered / wrapped / vertically centred) formatting
#    before D2's own style gets unified with the rest of the row.
# ---------------------------------------------------------------------
$ws2.Range("D2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

# Make every year header (D2:K2) match B2/C2's look (bold, centred, bordered)
$ws2.Range("B2").Copy()
$ws2.Range("D2:K2").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 3) New footnote style ("Sakkal Majalla", 12pt, vertically centred)
#    applied to the "Source : ..." line (A19) on both sheets
# ---------------------------------------------------------------------
foreach ($ws in @($ws1, $ws2)) {
    $src = $ws.Range("A19")
    $src.Font.Name = "Sakkal Majalla"
    $src.Font.Size = 12
    $src.VerticalAlignment = -4108   # xlCenter
}

# ---------------------------------------------------------------------
# 4) Views: make "Postes de Santé" the active/selected tab, leave A19
#    selected on "Centres de Santé", and zoom + select G5 on
#    "Postes de Santé"
# ---------------------------------------------------------------------
$ws1.Range("A19").Select()

$ws2.Activate()
$excel.ActiveWindow.Zoom = 97
$ws2.Range("G5").Select()

Write-Host "Workbook corrections applied"
